$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Lab 2: Introduction to Amazon EC2" -> append two trailing spaces plus a
#    bold / red / underlined "DONE" marker in its own run (mirrors the
#    existing "Lab 1 ... DONE" pattern already in the document).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Lab 2: Introduction to Amazon EC2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $markerStart = $rng.End
    $rng.InsertAfter("  DONE")
    $markerEnd = $rng.End

    $doneRange = $d.Range($markerStart + 2, $markerEnd)
    $doneRange.Font.Bold = $true
    $doneRange.Font.Color = 204
    $doneRange.Font.Underline = 1
}

# ---------------------------------------------------------------------------
# 2. "Lab 3: Introduction to Elastic Load Balancing" -> split after "Lab 3: "
#    and drop an (empty) bookmark named __DdeLink__90_1353281775 right at the
#    split point.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Lab 3: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Collapse(0)
    $bmRange = $d.Range($rng2.Start, $rng2.End)
    $d.Bookmarks.Add("__DdeLink__90_1353281775", $bmRange)
}
